# Fruta / hortaliza, semanal
# Weekly data refresh: insert a new "Ciruela" price record (Larry Ann,
# week of 2023-03-03) right after the existing row 66, pushing the
# subsequent records (previously rows 67-85) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 67:85 down to 68:86 by inserting a new blank row at 67.
$ws.Rows.Item(67).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A67").Value = 7
$ws.Range("B67").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C67").Value = "Ñuble"
$ws.Range("D67").Value = 44988
$ws.Range("E67").Value = 16
$ws.Range("F67").Value = "Fruta"
$ws.Range("G67").Value = 100103
$ws.Range("H67").Value = "Frutos de hueso (carozo)"
$ws.Range("I67").Value = 100103002
$ws.Range("J67").Value = "Ciruela"
$ws.Range("K67").Value = "Larry Ann"
$ws.Range("L67").Value = "Primera"
$ws.Range("M67").Value = 60
$ws.Range("N67").Value = 10000
$ws.Range("O67").Value = 10000
$ws.Range("P67").Value = 10000
$ws.Range("Q67").Value = "$/bandeja 18 kilos granel"
$ws.Range("R67").Value = "Región de O'Higgins"
$ws.Range("S67").Value = 556
$ws.Range("T67").Value = 18
